# Remove the trailing "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph
# and the "(c) 2020 . Contact: ..." footer paragraph that follow the
# bibliography entry for FLEMMING/GONCALVES, collapsing the blank paragraph
# that used to separate them from the rest of the document.

$d = $word.ActiveDocument

# Locate the start of the "Ver no Jupiter ..." paragraph.
$findStart = $d.Content
$okStart = $findStart.Find.Execute(
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

# Locate the end of the copyright/footer paragraph (search on an ASCII-only
# substring so we don't depend on how the copyright glyph round-trips).
$findEnd = $d.Content
$okEnd = $findEnd.Find.Execute(
    "Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if ($okStart -and $okEnd) {
    # Extend one character to the left so the blank paragraph mark right
    # before "Ver no Jupiter ..." is swallowed too (only one of the two
    # blank paragraphs that used to flank this block should survive), and
    # one character past the end of the copyright text so its own
    # paragraph mark is removed as well.
    $deleteRange = $d.Range($findStart.Start - 1, $findEnd.End + 1)
    $deleteRange.Delete()
}
